$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.709.27'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '3.091.15'
$ws.Range("E3").Value = '  -2.17%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'589.74"
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = "'153.62"
$ws.Range("E6").Value = '  +4.29%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = "'0.541"
$ws.Range("E8").Value = '  +1.82%  '
$ws.Range("D9").Value = '3.089.53'
$ws.Range("E9").Value = '  -2.05%  '
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("D14").Value = "'0.0000244"
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("D15").Value = '3.599.56'
$ws.Range("E15").Value = '  -2.24%  '
$ws.Range("E16").Value = '  -1.94%  '
$ws.Range("D17").Value = "'7.19"
$ws.Range("D18").Value = '63.666.22'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = '3.086.99'
$ws.Range("E19").Value = '  -2.27%  '
$ws.Range("D20").Value = "'472.32"
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Value = "'14.66"
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("D22").Value = "'0.729"
$ws.Range("E22").Value = '  -1.23%  '
$ws.Range("D23").Value = "'7.57"
$ws.Range("E23").Value = '  +0.57%  '
$ws.Range("D24").Value = "'2.40"
$ws.Range("E24").Value = '  +2.31%  '
$ws.Range("D25").Value = "'13.23"
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("D26").Value = "'81.67"
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = "'2.69"
$ws.Range("E29").Value = '  -1.39%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = "'7.30"
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("E31").Value = '  -1.47%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("E33").Value = '  +4.47%  '
$ws.Range("D34").Value = "'27.41"
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("D35").Value = '0.0₃0854'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").Value = "'3.43"
$ws.Range("E37").Value = '  +4.86%  '
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").Value = "'2.26"
$ws.Range("E39").Value = '  -4.40%  '
$ws.Range("E40").Value = '  +1.73%  '
$ws.Range("D41").Value = "'50.83"
$ws.Range("E41").Value = '  -1.89%  '
$ws.Range("D42").Value = "'449.73"
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("E43").Value = '  -2.42%  '
$ws.Range("E44").Value = '  -2.02%  '
$ws.Range("D45").Value = "'40.25"
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").Value = '2.836.61'
$ws.Range("E46").Value = '  -3.50%  '
$ws.Range("E47").Value = '  -0.66%  '
$ws.Range("D48").Value = "'130.16"
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("D49").Value = "'25.67"
$ws.Range("E49").Value = '  +4.78%  '
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("E51").Value = '  +0.45%  '
